$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)
$r = $ws2.Range("L1")

# The cell currently reads:
#   "Indemnités" (bold run) + "\n(Travail en équipes, le dimanche, de nuit et autres
#   primes de pénibilité, 1/12 de la somme annuelle)" (normal run)
# Remove the trailing ", 1/12 de la somme annuelle" from the second run, keeping the
# rest of the text and the rich-text run formatting (bold title / normal description).
$full = $r.Text
$toRemove = ", 1/12 de la somme annuelle"
$idx = $full.IndexOf($toRemove)
if ($idx -ge 0) {
    $c = $r.Characters($idx + 1, $toRemove.Length)
    $c.Text = ""
}

# Re-apply the run-level formatting, which the text splice above can otherwise disturb.
$newFull = $r.Text
$title = "Indemnités"
$titleLen = $title.Length
$restLen = $newFull.Length - $titleLen

$titleRun = $r.Characters(1, $titleLen)
$titleRun.Font.Bold = $true
$titleRun.Font.Size = 9
$titleRun.Font.Name = "Arial"

$descRun = $r.Characters($titleLen + 1, $restLen)
$descRun.Font.Bold = $false
$descRun.Font.Size = 9
$descRun.Font.Name = "Arial"
$descRun.Font.Color = 0

# The row height reflows to fit the (now shorter) wrapped text.
$ws2.Rows.Item(1).RowHeight = 259.5
